function Set-TextValue($Cell, $Text) {
    if ($Text -match '^-?\d+(\.\d+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "63.077.53"
Set-TextValue $ws.Range("E2") "  -0.27%  "
Set-TextValue $ws.Range("E3") "  -0.15%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "586.33"
Set-TextValue $ws.Range("E5") "  -0.59%  "
Set-TextValue $ws.Range("D6") "151.61"
Set-TextValue $ws.Range("E6") "  -0.99%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("E8") "  -1.92%  "
Set-TextValue $ws.Range("D9") "3.049.60"
Set-TextValue $ws.Range("E9") "  -0.46%  "
Set-TextValue $ws.Range("E10") "  -2.12%  "
Set-TextValue $ws.Range("D11") "5.86"
Set-TextValue $ws.Range("E11") "  +0.54%  "
Set-TextValue $ws.Range("D12") "0.448"
Set-TextValue $ws.Range("E12") "  -2.77%  "
Set-TextValue $ws.Range("D13") "0.0000235"
Set-TextValue $ws.Range("E13") "  -2.22%  "
Set-TextValue $ws.Range("D14") "36.25"
Set-TextValue $ws.Range("E14") "  -2.54%  "
Set-TextValue $ws.Range("E15") "  +1.98%  "
Set-TextValue $ws.Range("D16") "3.550.97"
Set-TextValue $ws.Range("E16") "  -0.37%  "
Set-TextValue $ws.Range("D17") "7.14"
Set-TextValue $ws.Range("E17") "  -0.95%  "
Set-TextValue $ws.Range("D18") "63.046.85"
Set-TextValue $ws.Range("E18") "  -0.37%  "
Set-TextValue $ws.Range("D19") "3.050.67"
Set-TextValue $ws.Range("E19") "  -0.48%  "
Set-TextValue $ws.Range("D20") "476.67"
Set-TextValue $ws.Range("E20") "  +0.22%  "
Set-TextValue $ws.Range("D21") "14.27"
Set-TextValue $ws.Range("E21") "  -2.57%  "
Set-TextValue $ws.Range("D22") "0.706"
Set-TextValue $ws.Range("E23") "  -0.28%  "
Set-TextValue $ws.Range("D24") "2.42"
Set-TextValue $ws.Range("E24") "  +1.18%  "
Set-TextValue $ws.Range("D25") "82.32"
Set-TextValue $ws.Range("E25") "  +1.51%  "
Set-TextValue $ws.Range("D26") "12.71"
Set-TextValue $ws.Range("E26") "  -1.62%  "
Set-TextValue $ws.Range("D27") "10.56"
Set-TextValue $ws.Range("E27") "  +5.48%  "
Set-TextValue $ws.Range("E29") "  +0.02%  "
Set-TextValue $ws.Range("D30") "2.67"
Set-TextValue $ws.Range("E30") "  -0.29%  "
Set-TextValue $ws.Range("E31") "  +0.20%  "
Set-TextValue $ws.Range("E32") "  +0.01%  "
Set-TextValue $ws.Range("D33") "27.66"
Set-TextValue $ws.Range("E33") "  +1.70%  "
Set-TextValue $ws.Range("E34") "  -2.45%  "
Set-TextValue $ws.Range("E35") "  +1.05%  "
Set-TextValue $ws.Range("D36") "0.0₃0819"
Set-TextValue $ws.Range("E36") "  -2.83%  "
Set-TextValue $ws.Range("E37") "  -2.63%  "
Set-TextValue $ws.Range("D38") "5.91"
Set-TextValue $ws.Range("E38") "  -3.10%  "
Set-TextValue $ws.Range("E39") "  +0.07%  "
Set-TextValue $ws.Range("D40") "9.23"
Set-TextValue $ws.Range("E40") "  -0.52%  "
Set-TextValue $ws.Range("E41") "  -0.05%  "
Set-TextValue $ws.Range("D42") "434.01"
Set-TextValue $ws.Range("E42") "  -2.24%  "
Set-TextValue $ws.Range("D43") "0.288"
Set-TextValue $ws.Range("E43") "  +0.76%  "
Set-TextValue $ws.Range("E44") "  +3.08%  "
Set-TextValue $ws.Range("D45") "0.0360"
Set-TextValue $ws.Range("E45") "  -0.49%  "
Set-TextValue $ws.Range("D46") "2.825.70"
Set-TextValue $ws.Range("E46") "  +0.99%  "
Set-TextValue $ws.Range("D47") "38.36"
Set-TextValue $ws.Range("E47") "  -4.42%  "
Set-TextValue $ws.Range("D48") "128.41"
Set-TextValue $ws.Range("E48") "  -2.55%  "
Set-TextValue $ws.Range("E49") "  -0.04%  "
Set-TextValue $ws.Range("D50") "25.12"
Set-TextValue $ws.Range("E50") "  +0.14%  "
Set-TextValue $ws.Range("E51") "  -1.64%  "
